# Fix issue with "Trenching" input to installation module.
# Replace incorrectly removed "ROV class" input.
# Add DateTimeDict for date outputs from installation module.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("rov")

# --- Insert a new "Class" column before the existing column C (Depth Rating) ---
$ws.Columns("C:C").Insert()
$ws.Columns("C:C").ColumnWidth = 17.983072916666668

# Capture the old "Name"-category values that used to live in column B
# (they describe the general Inspection class / Workclass category)
$b2 = $ws.Range("B2").Value2
$b3 = $ws.Range("B3").Value2
$b4 = $ws.Range("B4").Value2
$b5 = $ws.Range("B5").Value2
$b6 = $ws.Range("B6").Value2

# New column C holds the general category (what column B used to hold)
$ws.Range("C1").Value = "Class"
$ws.Range("C2").Value = $b2
$ws.Range("C3").Value = $b3
$ws.Range("C4").Value = $b4
$ws.Range("C5").Value = $b5
$ws.Range("C6").Value = $b6

# Column B now holds the specific ROV class name for each piece of equipment
$ws.Range("B2").Value = "Inspection class 1"
$ws.Range("B3").Value = "Inspection class 2"
$ws.Range("B4").Value = "Workclass 1"
$ws.Range("B5").Value = "Workclass 2"
$ws.Range("B6").Value = "Workclass 3"

# Fix up a couple of cell styles that need to match their row's formatting
$ws.Range("B2").Copy()
$ws.Range("B3").PasteSpecial(-4122)
$ws.Range("B4").Copy()
$ws.Range("B6").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Re-anchor the cell comments that shifted one column to the right ---
$commentRefs = @("N2", "P2", "P3", "G5", "N5", "O5", "P5")
$commentTexts = @{}
foreach ($r in $commentRefs) {
  $commentTexts[$r] = $ws.Range($r).Comment.Text()
}
foreach ($r in $commentRefs) {
  $ws.Range($r).Comment.Delete()
}
$ws.Range("O2").AddComment($commentTexts["N2"])
$ws.Range("Q2").AddComment($commentTexts["P2"])
$ws.Range("Q3").AddComment($commentTexts["P3"])
$ws.Range("H5").AddComment($commentTexts["G5"])
$ws.Range("O5").AddComment($commentTexts["N5"])
$ws.Range("P5").AddComment($commentTexts["O5"])
$ws.Range("Q5").AddComment($commentTexts["P5"])

# --- Make "rov" the active sheet/tab again (it had been left on cable_burial) ---
$ws.Activate()
$ws.Range("B8").Select()
